# Weekly meeting minutes update
# - F2: tweak wording
# - F3: clear the now-obsolete bullet
# - Row 5: fill in the 04.06 meeting entry (values + matching cell formatting)
# - Selection moves to E5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F2: reworded bullet --------------------------------------------------
$ws.Range("F2").Value2 = "-functii lamda, clase statice"

# --- F3: old bullet removed ------------------------------------------------
$ws.Range("F3").ClearContents()

# --- Row 5: new "04.06" meeting entry --------------------------------------
# Pull the cell formatting (fill/border/wrap) from equivalent cells first so
# the new row matches the rest of the table, then set the actual values.
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("E5").PasteSpecial(-4122)

$ws.Range("G2").Copy()
$ws.Range("G5").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

$ws.Range("A5").Value2 = "04.06"
$ws.Range("B5").Value2 = "clasa GPIO, exemplu friend"
$ws.Range("C5").Value2 = "exemplu friend, branch proiect, video"
$ws.Range("E5").Value2 = "branch, scheletul clasei, instantiere"
$ws.Range("G5").Value2 = "Andrei: CAN`nGabriel: USB`nAna: Timer`nAlina:RTC`nGPIO pe branch"

# Row height for the new, taller (multi-line) entry
$ws.Rows.Item(5).RowHeight = 75

# --- Selection moves to E5 --------------------------------------------------
$ws.Range("E5").Select()
